# Sliding window results: zero out the "IPC PO" prediction column (C),
# recompute DELTA (D) as -IPC RO and DELTA^2 (E) as IPC RO^2, then refresh
# the TOTAL (row 52) and MSE (row 53) summary cells to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 51

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2

    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = -$b
    $ws.Cells.Item($r, 5).Value = $b * $b
}

# TOTAL row: C52 sums the DELTA column, E52 sums the DELTA^2 column.
$deltaSum = 0
$delta2Sum = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $deltaSum += $ws.Cells.Item($r, 4).Value2
    $delta2Sum += $ws.Cells.Item($r, 5).Value2
}
$ws.Cells.Item(52, 3).Value = $deltaSum
$ws.Cells.Item(52, 5).Value = $delta2Sum

# MSE row: E53 is the mean of DELTA^2 over the 50 data points.
$ws.Cells.Item(53, 5).Value = $ws.Cells.Item(52, 5).Value2 / ($lastRow - $firstRow + 1)
